$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.485.54'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').Value = '3.174.49'
$ws.Range('E3').Value = '  -4.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.79'
$ws.Range('E5').Value = '  -2.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.32'
$ws.Range('E6').Value = '  -5.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '3.175.38'
$ws.Range('E8').Value = '  -4.18%  '
$ws.Range('E9').Value = '  -3.59%  '
$ws.Range('E10').Value = '  -5.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.28'
$ws.Range('E11').Value = '  -5.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.452'
$ws.Range('E12').Value = '  -4.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000235'
$ws.Range('E13').Value = '  -5.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.17'
$ws.Range('E14').Value = '  -4.81%  '
$ws.Range('D15').Value = '3.690.60'
$ws.Range('E15').Value = '  -4.42%  '
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('D17').Value = '3.165.66'
$ws.Range('E17').Value = '  -4.46%  '
$ws.Range('D18').Value = '62.429.38'
$ws.Range('E18').Value = '  -2.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.54'
$ws.Range('E19').Value = '  -5.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '455.77'
$ws.Range('E20').Value = '  -5.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.89'
$ws.Range('E21').Value = '  -2.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.704'
$ws.Range('E22').Value = '  -4.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.63'
$ws.Range('E23').Value = '  -5.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.42'
$ws.Range('E24').Value = '  -0.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.81'
$ws.Range('E25').Value = '  -2.43%  '
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.68'
$ws.Range('E27').Value = '  -3.17%  '
$ws.Range('B28').Value = 'FirstDigitalUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('E29').Value = '  -5.29%  '
$ws.Range('E30').Value = '  -4.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.02'
$ws.Range('E31').Value = '  -6.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.29'
$ws.Range('E32').Value = '  -7.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.104'
$ws.Range('E33').Value = '  -2.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.40'
$ws.Range('E34').Value = '  -6.16%  '
$ws.Range('E35').Value = '  -6.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.87'
$ws.Range('E36').Value = '  -2.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.17'
$ws.Range('E37').Value = '  -3.44%  '
$ws.Range('D38').Value = '0.0₃0696'
$ws.Range('E38').Value = '  -7.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0385'
$ws.Range('E39').Value = '  -4.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '413.68'
$ws.Range('E40').Value = '  -4.34%  '
$ws.Range('D41').Value = '2.904.14'
$ws.Range('E41').Value = '  -4.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.68'
$ws.Range('E42').Value = '  -2.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.00'
$ws.Range('E43').Value = '  -4.93%  '
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.15'
$ws.Range('E45').Value = '  -2.73%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.250'
$ws.Range('E46').Value = '  -6.32%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.15'
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '125.07'
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.38'
$ws.Range('E50').Value = '  -4.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.111'
$ws.Range('E51').Value = '  -3.69%  '
